$d = $word.ActiveDocument

# Step 0: locate the title run and capture its range boundaries.
$titleRange = $d.Content
$titleRange.Find.Execute("Escuela Secundaria Federal") | Out-Null
$startPos = $titleRange.Start
$textLen = $titleRange.End - $titleRange.Start

# Step 1: replace the plain title run with a 1-char placeholder run so we
# have a small, known-length anchor that sits *between* the existing
# w:bookmarkStart/w:bookmarkEnd pair (Word keeps the bookmark wrapped
# around whatever occupies that position).
$r1 = $d.Range($startPos, $startPos + $textLen)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:t>X</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# Step 2: replace that placeholder (still nested inside the bookmark) with
# the final formatted run containing the text box/drawing. Because the
# replaced range is strictly inside the bookmark span, InsertXML emits the
# new run *after* w:bookmarkEnd, matching the target ordering.
$r2 = $d.Range($startPos, $startPos + 1)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="6E2BC083" wp14:editId="1F854612"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>0</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>0</wp:posOffset></wp:positionV><wp:extent cx="1828800" cy="1828800"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapSquare wrapText="bothSides"/><wp:docPr id="1" name="Cuadro de texto 1"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1828800" cy="1828800"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln><a:effectLst/></wps:spPr><wps:txbx><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:color w:val="70AD47"/><w:spacing w:val="10"/><w:sz w:val="72"/><w:szCs w:val="72"/><w14:glow w14:rad="38100"><w14:schemeClr w14:val="accent1"><w14:alpha w14:val="60000"/></w14:schemeClr></w14:glow><w14:textOutline w14:w="9525" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="accent1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline><w14:textFill><w14:solidFill><w14:srgbClr w14:val="70AD47"><w14:tint w14:val="1000"/></w14:srgbClr></w14:solidFill></w14:textFill></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="70AD47"/><w:spacing w:val="10"/><w:sz w:val="72"/><w:szCs w:val="72"/><w14:glow w14:rad="38100"><w14:schemeClr w14:val="accent1"><w14:alpha w14:val="60000"/></w14:schemeClr></w14:glow><w14:textOutline w14:w="9525" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="accent1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline><w14:textFill><w14:solidFill><w14:srgbClr w14:val="70AD47"><w14:tint w14:val="1000"/></w14:srgbClr></w14:solidFill></w14:textFill></w:rPr><w:t>Escuela Secundaria Federal</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="none" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shapetype w14:anchorId="6E2BC083" id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="Cuadro de texto 1" o:spid="_x0000_s1026" type="#_x0000_t202" style="position:absolute;left:0;text-align:left;margin-left:0;margin-top:0;width:2in;height:2in;z-index:251659264;visibility:visible;mso-wrap-style:none;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQCFdGvFJwIAAFwEAAAOAAAAZHJzL2Uyb0RvYy54bWysVE2P2jAQvVfqf7B8LwFEWxoRVpQVVSW0&#10;uxJb7dk4DomUeCzbkNBf32cnsHTbU9WLM9+emfecxV3X1OykrKtIZ3wyGnOmtKS80oeM/3jefJhz&#10;5rzQuahJq4yfleN3y/fvFq1J1ZRKqnNlGYpol7Ym46X3Jk0SJ0vVCDciozScBdlGeKj2kORWtKje&#10;1Ml0PP6UtGRzY0kq52C97518GesXhZL+sSic8qzOOHrz8bTx3IczWS5EerDClJUc2hD/0EUjKo1L&#10;r6XuhRfsaKs/SjWVtOSo8CNJTUJFUUkVZ8A0k/GbaXalMCrOguU4c12T+39l5cPpybIqB3acadEA&#10;ovVR5JZYrphXnSc2CUtqjUsRuzOI9t1X6kLCYHcwhtm7wjbhi6kY/Fj3+bpiVGIyJM2n8/kYLgnf&#10;RUGd5DXdWOe/KWpYEDJugWFcrThtne9DLyHhNk2bqq5hF2mtfzOgZm9RkQhDdpik7zhIvtt3wxh7&#10;ys+YzlJPEmfkpkIHW+H8k7BgBboG0/0jjqKmNuM0SJyVZH/+zR7iARa8nLVgWcY1ngFn9XcNEL9M&#10;ZrNAyqjMPn6eQrG3nv2tRx+bNYHGAAq9RTHE+/oiFpaaFzyHVbgTLqElbs64v4hr3zMfz0mq1SoG&#10;gYZG+K3eGRlKhwWG7T53L8KaAYLAgwe6sFGkb5DoY0OmM6ujBx4RprDefqeANyigcAR6eG7hjdzq&#10;Mer1p7D8BQAA//8DAFBLAwQUAAYACAAAACEAS4kmzdYAAAAFAQAADwAAAGRycy9kb3ducmV2Lnht&#10;bEyP0U7DMAxF35H4h8hIvLF0FaBSmk5owDMw+ACvMU1p41RNthW+HoOQxovlq2tdn1utZj+oPU2x&#10;C2xguchAETfBdtwaeHt9vChAxYRscQhMBj4pwqo+PamwtOHAL7TfpFZJCMcSDbiUxlLr2DjyGBdh&#10;JBbvPUwek8ip1XbCg4T7QedZdq09diwfHI60dtT0m503UGT+qe9v8ufoL7+WV259Hx7GD2POz+a7&#10;W1CJ5nQ8hh98QYdamLZhxzaqwYAUSb9TvLwoRG7/Fl1X+j99/Q0AAP//AwBQSwECLQAUAAYACAAA&#10;ACEAtoM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQIt&#10;ABQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAAAAAAAC8BAABfcmVscy8ucmVsc1BLAQIt&#10;ABQABgAIAAAAIQCFdGvFJwIAAFwEAAAOAAAAAAAAAAAAAAAAAC4CAABkcnMvZTJvRG9jLnhtbFBL&#10;AQItABQABgAIAAAAIQBLiSbN1gAAAAUBAAAPAAAAAAAAAAAAAAAAAIEEAABkcnMvZG93bnJldi54&#10;bWxQSwUGAAAAAAQABADzAAAAhAUAAAAA&#10;" filled="f" stroked="f"><v:fill o:detectmouseclick="t"/><v:textbox style="mso-fit-shape-to-text:t"><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:color w:val="70AD47"/><w:spacing w:val="10"/><w:sz w:val="72"/><w:szCs w:val="72"/><w14:glow w14:rad="38100"><w14:schemeClr w14:val="accent1"><w14:alpha w14:val="60000"/></w14:schemeClr></w14:glow><w14:textOutline w14:w="9525" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="accent1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline><w14:textFill><w14:solidFill><w14:srgbClr w14:val="70AD47"><w14:tint w14:val="1000"/></w14:srgbClr></w14:solidFill></w14:textFill></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="70AD47"/><w:spacing w:val="10"/><w:sz w:val="72"/><w:szCs w:val="72"/><w14:glow w14:rad="38100"><w14:schemeClr w14:val="accent1"><w14:alpha w14:val="60000"/></w14:schemeClr></w14:glow><w14:textOutline w14:w="9525" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="accent1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline><w14:textFill><w14:solidFill><w14:srgbClr w14:val="70AD47"><w14:tint w14:val="1000"/></w14:srgbClr></w14:solidFill></w14:textFill></w:rPr><w:t>Escuela Secundaria Federal</w:t></w:r></w:p></w:txbxContent></v:textbox><w10:wrap type="square"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)

Write-Host "Title converted to WordArt text box."
